# Post AL test update - exemplo_encode.xlsx
# Fix object (lure and target) position: x-coordinate 250 -> 300
# lure_pos / target_pos values "(250, xxx)" become "(300, xxx)"
#
# Row2: target_pos (F2) was "(250, -200)" -> "(300, -200)"
#       lure_pos   (J2) was "(250, 200)"  -> "(300, 200)"
# Row3: target_pos (F3) was "(250, 200)"  -> "(300, 200)"
#       lure_pos   (J3) was "(250, -200)" -> "(300, -200)"
#
# Also move the active selection to F3 (was M4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the "(300, 200)" string first so it lands before "(300, -200)"
# in the shared-strings table, matching the new positions column order.
$ws.Range("J2").Value = "(300, 200)"
$ws.Range("F2").Value = "(300, -200)"
$ws.Range("F3").Value = "(300, 200)"
$ws.Range("J3").Value = "(300, -200)"

# Update the saved selection/active cell on the sheet.
$ws.Range("F3").Select()
